# Contrast Map.xlsx edit — move religion_muslim/religion_other "Intercaste" group
# to a new Sheet2, and insert new "both_general/both_obc/both_scst" caste
# contrast rows (for both h_dm/w_dm and h_htn/w_htn) into Sheet1 in their place.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Add Sheet2 right after Sheet1
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ------------------------------------------------------------------
# 2. Populate Sheet2 with the rows that used to live in Sheet1
#    (the religion_muslim / religion_other contrast rows for h_dm/w_dm
#    and h_htn/w_htn), now labelled "W6 R2"/"H6 R2"/"W6 R3"/"H6 R3".
# ------------------------------------------------------------------
$sheet2Data = @(
        @("h_dm at religion_muslim=0", "W6 R2", "Religion: Hindu"),
        @("h_dm at religion_muslim=1", "W6 R2", "Religion: Muslim"),
        @("Interaction of h_dm:religion_muslim", "W6 R2", $null),
        @("w_dm at religion_muslim=0", "H6 R2", "Religion: Hindu"),
        @("w_dm at religion_muslim=1", "H6 R2", "Religion: Muslim"),
        @("Interaction of w_dm:religion_muslim", "H6 R2", $null),
        @("h_dm at religion_other=0", "W6 R3", $null),
        @("h_dm at religion_other=1", "W6 R3", "Religion: Other"),
        @("Interaction of h_dm:religion_other", "W6 R3", $null),
        @("w_dm at religion_other=0", "H6 R3", $null),
        @("w_dm at religion_other=1", "H6 R3", "Religion: Other"),
        @("Interaction of w_dm:religion_other", "H6 R3", $null),
        @("h_htn at religion_muslim=0", "W6 R2", "Religion: Hindu"),
        @("h_htn at religion_muslim=1", "W6 R2", "Religion: Muslim"),
        @("Interaction of h_htn:religion_muslim", "W6 R2", $null),
        @("w_htn at religion_muslim=0", "H6 R2", "Religion: Hindu"),
        @("w_htn at religion_muslim=1", "H6 R2", "Religion: Muslim"),
        @("Interaction of w_htn:religion_muslim", "H6 R2", $null),
        @("h_htn at religion_other=0", "W6 R3", $null),
        @("h_htn at religion_other=1", "W6 R3", "Religion: Other"),
        @("Interaction of h_htn:religion_other", "W6 R3", $null),
        @("w_htn at religion_other=0", "H6 R3", $null),
        @("w_htn at religion_other=1", "H6 R3", "Religion: Other"),
        @("Interaction of w_htn:religion_other", "H6 R3", $null)
)

$r = 1
foreach ($row in $sheet2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) { $ws2.Cells.Item($r, 3).Value = $row[2] }
    $r++
}

# ------------------------------------------------------------------
# 3. Remove those rows from Sheet1 (delete bottom block first so the
#    row numbers of the top block don't shift before we delete it).
# ------------------------------------------------------------------
$ws1.Rows("122:133").Delete()
$ws1.Rows("56:67").Delete()

# ------------------------------------------------------------------
# 4. Make room in Sheet1 for the new caste contrast rows:
#    - 18 rows at 56:73   -> h_dm/w_dm at both_general/both_obc/both_scst
#    - 18 rows at 128:145 -> h_htn/w_htn at both_general/both_obc/both_scst
# ------------------------------------------------------------------
$ws1.Rows("56:73").Insert()
$ws1.Rows("128:145").Insert()

$dmData = @(
        @("h_dm at both_general=0", "W6 General", "Caste: Intercaste"),
        @("h_dm at both_general=1", "W6 General", "Caste: Both General"),
        @("Interaction of h_dm:both_general", "W6 General", $null),
        @("w_dm at both_general=0", "H6 General", "Caste: Intercaste"),
        @("w_dm at both_general=1", "H6 General", "Caste: Both General"),
        @("Interaction of w_dm:both_general", "H6 General", $null),
        @("h_dm at both_obc=0", "W6 OBC", $null),
        @("h_dm at both_obc=1", "W6 OBC", "Caste: Both OBC"),
        @("Interaction of h_dm:both_obc", "W6 OBC", $null),
        @("w_dm at both_obc=0", "H6 OBC", $null),
        @("w_dm at both_obc=1", "H6 OBC", "Caste: Both OBC"),
        @("Interaction of w_dm:both_obc", "H6 OBC", $null),
        @("h_dm at both_scst=0", "W6 SCST", $null),
        @("h_dm at both_scst=1", "W6 SCST", "Caste: Both SCST"),
        @("Interaction of h_dm:both_scst", "W6 SCST", $null),
        @("w_dm at both_scst=0", "H6 SCST", $null),
        @("w_dm at both_scst=1", "H6 SCST", "Caste: Both SCST"),
        @("Interaction of w_dm:both_scst", "H6 SCST", $null)
)

$r = 56
foreach ($row in $dmData) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) { $ws1.Cells.Item($r, 3).Value = $row[2] }
    $r++
}

$htnData = @(
        @("h_htn at both_general=0", "W6 General", "Caste: Intercaste"),
        @("h_htn at both_general=1", "W6 General", "Caste: Both General"),
        @("Interaction of h_htn:both_general", "W6 General", $null),
        @("w_htn at both_general=0", "H6 General", "Caste: Intercaste"),
        @("w_htn at both_general=1", "H6 General", "Caste: Both General"),
        @("Interaction of w_htn:both_general", "H6 General", $null),
        @("h_htn at both_obc=0", "W6 OBC", $null),
        @("h_htn at both_obc=1", "W6 OBC", "Caste: Both OBC"),
        @("Interaction of h_htn:both_obc", "W6 OBC", $null),
        @("w_htn at both_obc=0", "H6 OBC", $null),
        @("w_htn at both_obc=1", "H6 OBC", "Caste: Both OBC"),
        @("Interaction of w_htn:both_obc", "H6 OBC", $null),
        @("h_htn at both_scst=0", "W6 SCST", $null),
        @("h_htn at both_scst=1", "W6 SCST", "Caste: Both SCST"),
        @("Interaction of h_htn:both_scst", "W6 SCST", $null),
        @("w_htn at both_scst=0", "H6 SCST", $null),
        @("w_htn at both_scst=1", "H6 SCST", "Caste: Both SCST"),
        @("Interaction of w_htn:both_scst", "H6 SCST", $null)
)

$r = 128
foreach ($row in $htnData) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) { $ws1.Cells.Item($r, 3).Value = $row[2] }
    $r++
}

# ------------------------------------------------------------------
# 5. Cosmetics: column B width on Sheet1, sheet view selections.
# ------------------------------------------------------------------
$ws1.Columns("B").ColumnWidth = 11.1796875

$ws1.Application.ActiveWindow.ScrollRow = 4
[void]$ws1.Range("C56:C73").Select()

[void]$ws2.Range("A1:XFD12").Select()

Write-Output "done"
